{"js": "async (context) => {\n  // Mapping of old math problem text -> new math problem text.\n  // Every value below is unique within the document, so a simple\n  // search-and-replace for each pair unambiguously targets the right run.\n  const replacements = [\n    [\"56\u00d779=4424\", \"27\u00d783=2241\"],\n    [\"15\u00d718=270\", \"79\u00d792=7268\"],\n    [\"91\u00d742=3822\", \"33\u00d722=726\"],\n    [\"42\u00d744=1848\", \"59\u00d797=5723\"],\n    [\"53\u00d799=5247\", \"58\u00d711=638\"],\n    [\"17\u00d762=1054\", \"40\u00d742=1680\"],\n    [\"76\u00d792=6992\", \"35\u00d749=1715\"],\n    [\"73\u00d798=7154\", \"69\u00d730=2070\"],\n    [\"16\u00d779=1264\", \"96\u00d795=9120\"],\n    [\"25\u00d744=1100\", \"83\u00d771=5893\"],\n    [\"57\u00d778=4446\", \"45\u00d749=2205\"],\n    [\"57\u00d715=855\", \"96\u00d780=7680\"],\n    [\"70\u00d768=4760\", \"92\u00d786=7912\"],\n    [\"46\u00d740=1840\", \"16\u00d755=880\"],\n    [\"12\u00d791=1092\", \"38\u00d739=1482\"],\n    [\"12\u00d743=516\", \"20\u00d772=1440\"],\n    [\"92\u00d773=6716\", \"48\u00d765=3120\"],\n    [\"30\u00d745=1350\", \"53\u00d796=5088\"],\n    [\"24\u00d789=2136\", \"13\u00d718=234\"],\n    [\"73\u00d752=3796\", \"52\u00d788=4576\"],\n    [\"35\u00d759=2065\", \"94\u00d714=1316\"],\n    [\"80\u00d711=880\", \"66\u00d792=6072\"],\n    [\"57\u00d763=3591\", \"32\u00d719=608\"],\n    [\"31\u00d755=1705\", \"67\u00d727=1809\"],\n    [\"79\u00d729=2291\", \"14\u00d723=322\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const r of results.items) {\n      r.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n  }\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# Mapping of old math problem text -> new math problem text.\n# Every value below is unique within the document, so Find/Replace for each\n# pair unambiguously targets the right run.\n$replacements = @(\n    @(\"56\u00d779=4424\", \"27\u00d783=2241\"),\n    @(\"15\u00d718=270\", \"79\u00d792=7268\"),\n    @(\"91\u00d742=3822\", \"33\u00d722=726\"),\n    @(\"42\u00d744=1848\", \"59\u00d797=5723\"),\n    @(\"53\u00d799=5247\", \"58\u00d711=638\"),\n    @(\"17\u00d762=1054\", \"40\u00d742=1680\"),\n    @(\"76\u00d792=6992\", \"35\u00d749=1715\"),\n    @(\"73\u00d798=7154\", \"69\u00d730=2070\"),\n    @(\"16\u00d779=1264\", \"96\u00d795=9120\"),\n    @(\"25\u00d744=1100\", \"83\u00d771=5893\"),\n    @(\"57\u00d778=4446\", \"45\u00d749=2205\"),\n    @(\"57\u00d715=855\", \"96\u00d780=7680\"),\n    @(\"70\u00d768=4760\", \"92\u00d786=7912\"),\n    @(\"46\u00d740=1840\", \"16\u00d755=880\"),\n    @(\"12\u00d791=1092\", \"38\u00d739=1482\"),\n    @(\"12\u00d743=516\", \"20\u00d772=1440\"),\n    @(\"92\u00d773=6716\", \"48\u00d765=3120\"),\n    @(\"30\u00d745=1350\", \"53\u00d796=5088\"),\n    @(\"24\u00d789=2136\", \"13\u00d718=234\"),\n    @(\"73\u00d752=3796\", \"52\u00d788=4576\"),\n    @(\"35\u00d759=2065\", \"94\u00d714=1316\"),\n    @(\"80\u00d711=880\", \"66\u00d792=6072\"),\n    @(\"57\u00d763=3591\", \"32\u00d719=608\"),\n    @(\"31\u00d755=1705\", \"67\u00d727=1809\"),\n    @(\"79\u00d729=2291\", \"14\u00d723=322\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdReplaceAll = 2 (there is exactly one occurrence of each $old, but\n    # ReplaceAll is used so the call is idempotent/safe either way)\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
